# Update 4/2/2021 - Alpha stage sphinx documentation - Unit test works
#
# Changes applied:
#  1. "Connectivity" sheet: just move the selection (no data change).
#  2. "Member" sheet: insert a new column C "Trans tag" holding the
#     OpenSees geometric-transformation tag used by each member, then
#     move the selection.
#  3. Add a brand new "Member transformation" worksheet (as the last /
#     4th tab) describing the two transformation vectors referenced by
#     the new "Trans tag" column, and leave it as the active sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Connectivity: move selection only ------------------------------
$wsConn = $wb.Worksheets.Item("Connectivity")
$wsConn.Range("E27").Select() | Out-Null

# --- 2. Member: insert "Trans tag" column -------------------------------
$wsMember = $wb.Worksheets.Item("Member")

# New column pushes A(m^2)...Ay(m^2) one slot to the right (C -> D, etc.)
$wsMember.Columns.Item(3).Insert()
$wsMember.Columns.Item(3).ColumnWidth = $wsMember.Columns.Item(2).ColumnWidth

$wsMember.Cells.Item(1, 3).Value = "Trans tag"
$wsMember.Cells.Item(2, 3).Value = 1
$wsMember.Cells.Item(3, 3).Value = 1
$wsMember.Cells.Item(4, 3).Value = 1
$wsMember.Cells.Item(5, 3).Value = 2
$wsMember.Cells.Item(6, 3).Value = 2

$wsMember.Range("E16").Select() | Out-Null

# --- 3. New "Member transformation" sheet --------------------------------
$wsTrans = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsTrans.Name = "Member transformation"

$wsTrans.Columns.Item(1).ColumnWidth = 15.33203125
$wsTrans.Columns.Item(2).ColumnWidth = 36.44140625
$wsTrans.Columns.Item(3).ColumnWidth = 8.88671875

$wsTrans.Cells.Item(1, 1).Value = "Transform tag "
$wsTrans.Cells.Item(1, 2).Value = "Name of variable saved in grillage class"
$wsTrans.Cells.Item(1, 3).Value = "Vector"

$wsTrans.Cells.Item(2, 1).Value = 1
$wsTrans.Cells.Item(2, 2).Value = "long"
$wsTrans.Cells.Item(2, 3).Value = "[0,0,1]"

$wsTrans.Cells.Item(3, 1).Value = 2
$wsTrans.Cells.Item(3, 2).Value = "trans"
$wsTrans.Cells.Item(3, 3).Value = "[1,0,0]"

$wsTrans.Range("C9").Select() | Out-Null
